$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name to reflect new "through" date
$ws.Name = "Through 2022-10-28"

# Update header label cell (I1) which holds the shared string "2022 (through 10-27)"
$ws.Range("I1").Value = "2022 (through 10-28)"

# Update October row's 2022 value (row 11, column I)
$ws.Range("I11").Value = 101

# Update Total row's 2022 value (row 14, column I)
$ws.Range("I14").Value = 1377
